$wb = $excel.ActiveWorkbook
$events = $wb.Worksheets.Item("Events")
$sources = $wb.Worksheets.Item("Sources")

# --- Events sheet -----------------------------------------------------
# Remove the SecurityRoomBG row (row 9) - it is being cut / folded out of scope.
$events.Rows.Item(9).Delete()

# Remove the CollectGlasses row (originally row 15, now row 14 after the
# previous delete) - superseded by the GlassesDrop entry.
$events.Rows.Item(14).Delete()

# Rename the movement-interface click event to match the implemented
# asset name.
$events.Cells.Item(7, 1).Value = "GroundClick"

# All remaining events have now been implemented in-engine.
for ($r = 2; $r -le 16; $r++) {
    $events.Cells.Item($r, 5).Value = "Implemented"
}

# --- Sources sheet ------------------------------------------------------
# New asset sourced for the coffee-grinding sound effect.
$sources.Cells.Item(18, 1).Value = "Grinding Beans for Cold Brew Coffee"
$sources.Hyperlinks.Add($sources.Cells.Item(18, 2), "https://freesound.org/people/munyeca/sounds/348148/")

$sources.Range("A18").Select()

# Leave focus back on the Events tab (it is the active sheet in the saved
# workbook), with the last-touched cell selected.
$events.Activate()
$events.Range("E16").Select()
